$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the header formatting from the last existing header cell (Z1) onto the
# two new header cells before setting their text.
$ws.Range("Z1").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AA1").Value = "Exp Constant"
$ws.Range("AB1").Value = "Exp Constant [dB]"

# New data values for rows 2-5
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 27).Value = 385250961.9682089
    $ws.Cells.Item($r, 28).Value = 85.85743731821252
}
